# Add the "HN vs LN" contrast row (row 23) to each of the four worksheets.

$wb = $excel.ActiveWorkbook

$data = @{
    "MaxExG"      = @(33177.352366666702, 16404.804152434699, 24, 2.0224168516966201, [double]"5.4413631590805603E-2")
    "SDExG"       = @(46.364920073333401, 11.426370763451001, 24, 4.05771185209906, [double]"4.5536108812006202E-4")
    "LocExG"      = @(27.686784930000101, 7.81636777870711, 24, 3.5421548363451798, [double]"1.65971444217679E-3")
    "MaxSlopeExG" = @(-196.98121933333201, 673.90769491605897, 24, -0.29229703239679899, 0.77257149993784702)
}

$highlighted = @("SDExG", "LocExG")

# Order matters only for which sheet/selection ends up "active" at the end,
# so process in tab order and finish on MaxSlopeExG (the originally active tab).
$sheetOrder = @("MaxExG", "SDExG", "LocExG", "MaxSlopeExG")

foreach ($sheetName in $sheetOrder) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $data[$sheetName]

    $ws.Range("A23").Value = "HN vs LN"
    $ws.Range("B23").Value = $values[0]
    $ws.Range("C23").Value = $values[1]
    $ws.Range("D23").Value = $values[2]
    $ws.Range("E23").Value = $values[3]
    $ws.Range("F23").Value = $values[4]

    if ($highlighted -contains $sheetName) {
        $ws.Range("F23").Interior.ColorIndex = 6
        $ws.Range("F23").Select() | Out-Null
    } else {
        $ws.Range("A23:F23").Select() | Out-Null
    }
}
